$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Turn the final (empty) paragraph into the bold "Architecture"
#    heading paragraph.
# ---------------------------------------------------------------
$archPara = $d.Paragraphs.Last
$archXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Architecture</w:t></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$archPara.Range.InsertXML($archXml)

# ---------------------------------------------------------------
# 2) Add a new bold "BillRive-Service" paragraph right after it
#    (BillRive is wrapped with spellcheck proofErr markers, same as
#    the rest of the document).
# ---------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null
$svcPara = $d.Paragraphs.Last
$svcXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>BillRive</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>-Service</w:t></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$svcPara.Range.InsertXML($svcXml)

# ---------------------------------------------------------------
# 3) Add the new bulleted list paragraph describing the Spring MVC
#    architecture video, using a brand-new list definition (numId 3).
# ---------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null
$listPara = $d.Paragraphs.Last
$listXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Spring MVC (REST) architecture. REST architecture following the guidelines from this video (</w:t></w:r><w:r><w:t>http://www.youtube.com/watch?v=5WXYw4J4QOU</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$listPara.Range.InsertXML($listXml)

# Give that paragraph a real bullet-list numbering definition (numId 3)
# so numbering.xml gains the matching abstractNum/num entries.
$lg = $word.ListGalleries.Item(1)
$lt = $lg.ListTemplates.Item(1)
$listPara.Range.ListFormat.ApplyListTemplate($lt)

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
